$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price/Volume columns so numeric-looking strings
# (e.g. "250.77") are kept as literal text instead of being parsed as numbers,
# matching the inline-string cells used by the source data feed.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '42.107.40'
$ws.Range('E2').Value = '  +1.62%  '
$ws.Range('D3').Value = '2.213.31'
$ws.Range('E3').Value = '  +0.85%  '
$ws.Range('D5').Value = '250.77'
$ws.Range('E5').Value = '  -1.13%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  -0.69%  '
$ws.Range('D7').Value = '68.09'
$ws.Range('E7').Value = '  -0.66%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '0.620'
$ws.Range('E9').Value = '  +6.23%  '
$ws.Range('D10').Value = '39.22'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').Value = '59.36'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').Value = '0.0937'
$ws.Range('E12').Value = '  -0.58%  '
$ws.Range('D13').Value = '7.07'
$ws.Range('E13').Value = '  -1.21%  '
$ws.Range('D14').Value = '0.104'
$ws.Range('E14').Value = '  -0.53%  '
$ws.Range('D15').Value = '2.548.71'
$ws.Range('E15').Value = '  +1.22%  '
$ws.Range('D16').Value = '0.868'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '14.48'
$ws.Range('E17').Value = '  -0.52%  '
$ws.Range('D18').Value = '2.210.71'
$ws.Range('E18').Value = '  +1.15%  '
$ws.Range('D19').Value = '42.005.17'
$ws.Range('E19').Value = '  +1.75%  '
$ws.Range('D20').Value = '0.0₃0960'
$ws.Range('E20').Value = '  +0.27%  '
$ws.Range('D21').Value = '72.35'
$ws.Range('E21').Value = '  +0.19%  '
$ws.Range('E22').Value = '  -2.07%  '
$ws.Range('D23').Value = '231.35'
$ws.Range('E23').Value = '  -0.92%  '
$ws.Range('E24').Value = '  -1.23%  '
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -5.57%  '
$ws.Range('E28').Value = '  -4.47%  '
$ws.Range('D29').Value = '3.68'
$ws.Range('E29').Value = '  -1.76%  '
$ws.Range('E30').Value = '  +1.17%  '
$ws.Range('D31').Value = '166.47'
$ws.Range('E31').Value = '  -1.95%  '
$ws.Range('D32').Value = '20.44'
$ws.Range('E32').Value = '  -1.26%  '
$ws.Range('D33').Value = '6.07'
$ws.Range('E33').Value = '  +10.18%  '
$ws.Range('D34').Value = '0.121'
$ws.Range('E34').Value = '  +1.70%  '
$ws.Range('D35').Value = '0.0778'
$ws.Range('E35').Value = '  +6.07%  '
$ws.Range('D36').Value = '0.122'
$ws.Range('E36').Value = '  -1.23%  '
$ws.Range('D37').Value = '26.51'
$ws.Range('E37').Value = '  -0.66%  '
$ws.Range('D38').Value = '4.59'
$ws.Range('E38').Value = '  -1.08%  '
$ws.Range('E39').Value = '  -1.05%  '
$ws.Range('E40').Value = '  +4.22%  '
$ws.Range('D41').Value = '2.23'
$ws.Range('E41').Value = '  +0.15%  '
$ws.Range('E42').Value = '  -0.80%  '
$ws.Range('D43').Value = '12.04'
$ws.Range('E43').Value = '  -6.36%  '
$ws.Range('D44').Value = '5.09'
$ws.Range('E44').Value = '  +2.87%  '
$ws.Range('D45').Value = '61.55'
$ws.Range('E45').Value = '  -4.97%  '
$ws.Range('E46').Value = '  -3.76%  '
$ws.Range('D47').Value = '8.57'
$ws.Range('E47').Value = '  -1.45%  '
$ws.Range('D48').Value = '0.100'
$ws.Range('E48').Value = '  -1.32%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').Value = '1.14'
$ws.Range('E50').Value = '  +0.03%  '
$ws.Range('E51').Value = '  +0.86%  '
